$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 765
$ws.Range("J19").Value = 812.2857
$ws.Range("L19").Value = 812.2857
$ws.Range("N19").Value = -1162.2857

$ws.Range("H41").Value = 591.1429000000001
$ws.Range("I41").Value = 467
$ws.Range("J41").Value = 684.25
$ws.Range("K41").Value = 467
$ws.Range("L41").Value = 684.25
$ws.Range("M41").Value = -27
$ws.Range("N41").Value = -1564.25

$ws.Range("H70").Value = 1994.4445
$ws.Range("I70").Value = 1850
$ws.Range("K70").Value = 5550
$ws.Range("M70").Value = -5280

$ws.Range("H73").Value = 1994.4445
$ws.Range("I73").Value = 1850
$ws.Range("K73").Value = 5550
$ws.Range("M73").Value = -4614

$ws.Range("H80").Value = 727.1667
$ws.Range("J80").Value = 1266.6666
$ws.Range("L80").Value = 3799.9998
$ws.Range("N80").Value = -5795.9998

$ws.Range("H83").Value = 727.1667
$ws.Range("J83").Value = 1266.6666
$ws.Range("L83").Value = 11399.9994
$ws.Range("N83").Value = -21383.9994

$ws.Range("H96").Value = 1306.25
$ws.Range("I96").Value = 1058.579
$ws.Range("K96").Value = 3175.737
$ws.Range("M96").Value = -1802.737

$ws.Range("H98").Value = 2485.7441
$ws.Range("I98").Value = 2204.8206
$ws.Range("K98").Value = 2204.8206
$ws.Range("M98").Value = -706.8206

$ws.Range("H103").Value = 566.25
$ws.Range("J103").Value = 420.6
$ws.Range("L103").Value = 1261.8
$ws.Range("N103").Value = -2433.8

$ws.Range("H122").Value = 2485.7441
$ws.Range("I122").Value = 2204.8206
$ws.Range("K122").Value = 6614.4618
$ws.Range("M122").Value = -4164.4618

$ws.Range("H131").Value = 556055.3
$ws.Range("I131").Value = 556055.3
$ws.Range("K131").Value = 1668165.9
$ws.Range("M131").Value = -1663125.9

$ws.Range("H138").Value = 2078.697
$ws.Range("I138").Value = 1426.4324
$ws.Range("J138").Value = 2467.9517
$ws.Range("K138").Value = 4279.2972
$ws.Range("L138").Value = 7403.855100000001
$ws.Range("M138").Value = 860.7028
$ws.Range("N138").Value = -17683.8551

$ws.Range("H141").Value = 3316.1333
$ws.Range("I141").Value = 3410.1428
$ws.Range("K141").Value = 10230.4284
$ws.Range("M141").Value = -5050.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 99468.914
$ws.Range("I61").Value = 2255.6365
$ws.Range("J61").Value = 181726.31
$ws.Range("K61").Value = 2255.6365
$ws.Range("L61").Value = 181726.31
$ws.Range("M61").Value = -2043.6365
$ws.Range("N61").Value = -182150.31

$ws.Range("H63").Value = 3219.8
$ws.Range("I63").Value = 3299.75
$ws.Range("K63").Value = 3299.75
$ws.Range("M63").Value = -2613.75

$ws.Range("H66").Value = 3219.8
$ws.Range("I66").Value = 3299.75
$ws.Range("K66").Value = 16498.75
$ws.Range("M66").Value = -13066.75

$ws.Range("H74").Value = 11315.279
$ws.Range("I74").Value = 1597.7742
$ws.Range("J74").Value = 36418.832
$ws.Range("K74").Value = 1597.7742
$ws.Range("L74").Value = 36418.832
$ws.Range("M74").Value = -723.7742000000001
$ws.Range("N74").Value = -38166.832

$ws.Range("H77").Value = 11315.279
$ws.Range("I77").Value = 1597.7742
$ws.Range("J77").Value = 36418.832
$ws.Range("K77").Value = 7988.871
$ws.Range("L77").Value = 182094.16
$ws.Range("M77").Value = -3620.871
$ws.Range("N77").Value = -190830.16

$ws.Range("H88").Value = 2275.7144
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 2000
$ws.Range("M88").Value = -1594

$ws.Range("H91").Value = 2275.7144
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 2000
$ws.Range("M91").Value = -596

$ws.Range("H97").Value = 1007.5
$ws.Range("I97").Value = 1011.6667
$ws.Range("J97").Value = 995
$ws.Range("K97").Value = 1011.6667
$ws.Range("L97").Value = 995
$ws.Range("M97").Value = -515.6667
$ws.Range("N97").Value = -1987

$ws.Range("H136").Value = 99468.914
$ws.Range("I136").Value = 2255.6365
$ws.Range("J136").Value = 181726.31
$ws.Range("K136").Value = 6766.9095
$ws.Range("L136").Value = 545178.9299999999
$ws.Range("M136").Value = -4216.9095
$ws.Range("N136").Value = -550278.9299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10436181
$ws.Range("I20").Value = 22229822
$ws.Range("J20").Value = 30026.412
$ws.Range("K20").Value = 22229822
$ws.Range("L20").Value = 30026.412
$ws.Range("M20").Value = -22229575
$ws.Range("N20").Value = -30520.412

$ws.Range("H86").Value = 38463270
$ws.Range("I86").Value = 1511.6666
$ws.Range("K86").Value = 1511.6666
$ws.Range("M86").Value = -388.6666

$ws.Range("H89").Value = 38463270
$ws.Range("I89").Value = 1511.6666
$ws.Range("K89").Value = 7558.333000000001
$ws.Range("M89").Value = -1942.333000000001

$ws.Range("H99").Value = 35951.54
$ws.Range("I99").Value = 41260.91
$ws.Range("K99").Value = 41260.91
$ws.Range("M99").Value = -39762.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2467.1667
$ws.Range("I19").Value = 2860.6
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 2860.6
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -2690.6
$ws.Range("N19").Value = -840

$ws.Range("H24").Value = 2467.1667
$ws.Range("I24").Value = 2860.6
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 2860.6
$ws.Range("L24").Value = 500
$ws.Range("M24").Value = -2690.6
$ws.Range("N24").Value = -840

$ws.Range("H31").Value = 10911.549
$ws.Range("J31").Value = 20143.875
$ws.Range("L31").Value = 20143.875
$ws.Range("N31").Value = -20733.875

$ws.Range("H34").Value = 10911.549
$ws.Range("J34").Value = 20143.875
$ws.Range("L34").Value = 20143.875
$ws.Range("N34").Value = -20547.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.53846
$ws.Range("I2").Value = 92.2
$ws.Range("J2").Value = 40.214287
$ws.Range("K2").Value = 553.2
$ws.Range("L2").Value = 241.285722
$ws.Range("M2").Value = -440.2
$ws.Range("N2").Value = -467.285722

$ws.Range("H7").Value = 185.46666
$ws.Range("J7").Value = 572.5
$ws.Range("L7").Value = 1717.5
$ws.Range("N7").Value = -1941.5

$ws.Range("H68").Value = 7499
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 7499
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 22497
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -24119

$ws.Range("H71").Value = 7499
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 7499
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 67491
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -75603

$ws.Range("H121").Value = 461.35715
$ws.Range("I121").Value = 387.18182
$ws.Range("K121").Value = 1161.54546
$ws.Range("M121").Value = 148.45454

$ws.Range("H131").Value = 1446.11
$ws.Range("I131").Value = 999
$ws.Range("J131").Value = 1455.2347
$ws.Range("K131").Value = 2997
$ws.Range("L131").Value = 4365.7041
$ws.Range("M131").Value = 2043
$ws.Range("N131").Value = -14445.7041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19499.75
$ws.Range("J70").Value = 26500
$ws.Range("L70").Value = 26500
$ws.Range("N70").Value = -27040

$ws.Range("H73").Value = 19499.75
$ws.Range("J73").Value = 26500
$ws.Range("L73").Value = 26500
$ws.Range("N73").Value = -28372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11769.195
$ws.Range("I136").Value = 9124.714
$ws.Range("J136").Value = 17465
$ws.Range("K136").Value = 27374.142
$ws.Range("L136").Value = 52395
$ws.Range("M136").Value = -24824.142
$ws.Range("N136").Value = -57495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1834.05
$ws.Range("I96").Value = 1855.7142
$ws.Range("J96").Value = 1822.3846
$ws.Range("K96").Value = 1855.7142
$ws.Range("L96").Value = 1822.3846
$ws.Range("M96").Value = -482.7141999999999
$ws.Range("N96").Value = -4568.3846

$ws.Range("H100").Value = 768.3
$ws.Range("J100").Value = 843.4
$ws.Range("L100").Value = 1686.8
$ws.Range("N100").Value = -2768.8
